# One additional historical fact from Jonathan Belcher added.
# A new row is inserted above the existing row 519 (1975 / "MBTA Red Line
# Harvard...") recording a new 1975 MBTA Orange Line fact; every row from
# the old 519 onward shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 519, pushing rows 519..655 down to 520..656.
$ws.Rows(519).Insert() | Out-Null

# Fill in the new fact.
$ws.Range("A519").Value = 1975
$ws.Range("B519").Value = "The relocated northern portion of the MBTA Orange line opened to Sullivan in April, to Wellington in September, and to Malden in December."

# Match the author's final cursor position/selection after typing the new
# text and pressing Enter (selection moves to the cell below, B520).
$ws.Range("B520").Select() | Out-Null
